# Generate Report for Handoff
#
# Updates the localization-status report to reflect a fresh handoff:
#  - "Handed back: in sync with en-US" -> "Ready for handoff" (Status cells)
#  - bump the handoff timestamps forward
#  - narrow the two "Status"-ish columns that used to be sized for the long
#    "Handed back: in sync with en-US" text, now that the text is shorter

$wb = $excel.ActiveWorkbook

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

# Target raw column width (as stored in the worksheet XML) is 17.2159881591797.
# Excel's COM ColumnWidth setter snaps to whole-pixel granularity (same
# pixel-grid rounding real Excel applies), so the nearest value reachable via
# ColumnWidth is 17.166666... ; feeding 16.3333... as the COM ColumnWidth
# lands exactly on that nearest reachable width.
$newWidth = 16.3333333333333

# --- Overview sheet --------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-21 11:03:26"

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-21 11:03:22"

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-21 11:03:26"

$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
